$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178289294242859
$ws.Range("B1").Value = 2.211082458496094
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.357163190841675
$ws.Range("E1").Value = 1.224979162216187
